# DataToValidateNameField.xlsx — rebuild the "Sheet1" test-data column.
# The previous fixture had 10 rows of miscellaneous name/validation samples
# (Sai, Konduru Bharath Sai, pavan789823, Telugu/Tamil/Hindi scripts, a
# super-long string, a huge number, special characters, ...) with ad-hoc
# per-cell formatting (wrapped text, a custom "JetBrains Mono" font, a wide
# column). The new fixture is a clean 7-row list of plain names used by the
# "invalidName" test.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old rows 8-10 (long Mahesh string w/ wrap formatting, the huge
# scientific-notation number, and the special-character string w/ custom
# font) entirely, rather than merely clearing their contents, so the used
# range shrinks back down with them.
$ws.Range("A8:A10").EntireRow.Delete()

# The surviving rows also get reshuffled/replaced with the new name list.
$ws.Range("A1").Value = "Bharath"
$ws.Range("A2").Value = "Pavan"
$ws.Range("A3").Value = "Dhruv"
$ws.Range("A4").Value = "Santosh"
$ws.Range("A5").Value = "Naveen"
$ws.Range("A6").Value = "Kondurur Bharath Sai"
$ws.Range("A7").Value = "KBS"

# Column A no longer needs the wide, manually-set width from before.
$ws.Columns.Item(1).ColumnWidth = 8.43

# Leave the selection on the last populated cell, like the saved file does.
$ws.Range("A7").Select() | Out-Null
